$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to text so numeric-looking values (e.g. "324.45")
# are stored as text, matching the original inlineStr cell type.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "27.706.62"
$ws.Range("E2").Value = "  -0.93%  "
$ws.Range("D3").Value = "1.756.94"
$ws.Range("E3").Value = "  -0.68%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "324.45"
$ws.Range("E5").Value = "  +0.87%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").Value = "0.4609"
$ws.Range("E7").Value = "  +8.02%  "
$ws.Range("D8").Value = "0.3602"
$ws.Range("E8").Value = "  -0.77%  "
$ws.Range("D9").Value = "0.07516"
$ws.Range("E9").Value = "  +0.53%  "
$ws.Range("D10").Value = "42.21"
$ws.Range("E10").Value = "  -2.78%  "
$ws.Range("E11").Value = "  +0.28%  "
$ws.Range("D13").Value = "20.77"
$ws.Range("E13").Value = "  -1.26%  "
$ws.Range("D14").Value = "6.021"
$ws.Range("E14").Value = "  -1.50%  "
$ws.Range("D15").Value = "7.122"
$ws.Range("E15").Value = "  -2.85%  "
$ws.Range("D16").Value = "1.760.51"
$ws.Range("E16").Value = "  -1.82%  "
$ws.Range("D17").Value = "92.32"
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("D18").Value = "0.00001067"
$ws.Range("E18").Value = "  +0.47%  "
$ws.Range("D19").Value = "0.06396"
$ws.Range("E19").Value = "  -0.15%  "
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("D21").Value = "16.80"
$ws.Range("E21").Value = "  -2.08%  "
$ws.Range("D22").Value = "5.804"
$ws.Range("E22").Value = "  -2.94%  "
$ws.Range("D23").Value = "27.761.91"
$ws.Range("E23").Value = "  -0.82%  "
$ws.Range("D24").Value = "11.27"
$ws.Range("E24").Value = "  -0.67%  "
$ws.Range("D25").Value = "2.105"
$ws.Range("E25").Value = "  -0.27%  "
$ws.Range("D26").Value = "164.15"
$ws.Range("E26").Value = "  +3.74%  "
$ws.Range("D27").Value = "20.37"
$ws.Range("E27").Value = "  +0.54%  "
$ws.Range("D28").Value = "1.960.16"
$ws.Range("E28").Value = "  -1.60%  "
$ws.Range("D29").Value = "2.076"
$ws.Range("E29").Value = "  -4.02%  "
$ws.Range("D30").Value = "126.44"
$ws.Range("E30").Value = "  +0.64%  "
$ws.Range("D31").Value = "1.055"
$ws.Range("E31").Value = "  -8.59%  "
$ws.Range("D32").Value = "0.09239"
$ws.Range("E32").Value = "  +3.57%  "
$ws.Range("D33").Value = "3.670"
$ws.Range("E33").Value = "  -1.51%  "
$ws.Range("D34").Value = "5.529"
$ws.Range("E34").Value = "  -1.86%  "
$ws.Range("D35").Value = "11.90"
$ws.Range("E35").Value = "  -4.58%  "
$ws.Range("D36").Value = "0.02302"
$ws.Range("E36").Value = "  -0.58%  "
$ws.Range("D37").Value = "0.2102"
$ws.Range("E37").Value = "  -0.42%  "
$ws.Range("D38").Value = "0.06034"
$ws.Range("E38").Value = "  +0.28%  "
$ws.Range("D39").Value = "0.6348"
$ws.Range("E39").Value = "  -0.23%  "
$ws.Range("D40").Value = "4.972"
$ws.Range("E40").Value = "  -1.44%  "
$ws.Range("E41").Value = "  +1.49%  "
$ws.Range("D42").Value = "1.380"
$ws.Range("E42").Value = "  -1.97%  "
$ws.Range("D43").Value = "7.808"
$ws.Range("E43").Value = "  -0.40%  "
$ws.Range("D44").Value = "13.30"
$ws.Range("E44").Value = "  -1.52%  "
$ws.Range("D45").Value = "0.5912"
$ws.Range("E45").Value = "  -0.40%  "
$ws.Range("D46").Value = "3.716"
$ws.Range("E46").Value = "  +0.52%  "
$ws.Range("D47").Value = "123.47"
$ws.Range("E47").Value = "  +0.58%  "
$ws.Range("D48").Value = "1.952"
$ws.Range("E48").Value = "  -2.91%  "
$ws.Range("D49").Value = "1.149"
$ws.Range("E49").Value = "  -3.66%  "
$ws.Range("D50").Value = "0.06869"
$ws.Range("E50").Value = "  +0.18%  "
$ws.Range("D51").Value = "72.31"
$ws.Range("E51").Value = "  -2.84%  "

# Restore default (Normal) style on the Price column so no stray
# number-format style index is left behind on the cells.
$priceRange.Style = "Normal"
